# Update plots for each sample
# - peak_table: last marker's (CYP2D6_XN, sample S2) measured m_height dropped to 800
# - allele_table: the corresponding reverse-direction "A" peak for CYP2D6_XN/S2 is now
#   detected (row 35), so its detection/measurement columns get populated
# - marker_table: CYP2D6_XN/S2 genotype now reads heterozygous (GA) instead of wildtype (GG)
# - genotype_result: overall sample genotype call updated accordingly

$wb = $excel.ActiveWorkbook

$peakTable = $wb.Worksheets.Item("peak_table")
$peakTable.Range("O18").Value = 800

$alleleTable = $wb.Worksheets.Item("allele_table")
$alleleTable.Range("K35").Value = 800
$alleleTable.Range("M35").Value = $true
$alleleTable.Range("N35").Value = 22
$alleleTable.Range("O35").Value = 72.05
$alleleTable.Range("P35").Value = 858
$alleleTable.Range("Q35").Value = "ok"
$alleleTable.Range("R35").Value = ""

$markerTable = $wb.Worksheets.Item("marker_table")
$markerTable.Range("G18").Value = "GA"
$markerTable.Range("H18").Value = "heterozygous"

$genotypeResult = $wb.Worksheets.Item("genotype_result")
$genotypeResult.Range("B2").Value = "*1/*10BX2|*1XN/*10B"
